$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 525, pushing existing rows 525:576 down to 526:577
$ws.Rows("525:525").Insert()

# Populate the newly inserted row 525 with the new data entry
$ws.Range("A525").Value = 8
$ws.Range("B525").Value = "Terminal La Palmera de La Serena"
$ws.Range("C525").Value = "Coquimbo"
$ws.Range("D525").Value = 45166
$ws.Range("E525").Value = 4
$ws.Range("F525").Value = 100112032
$ws.Range("G525").Value = "Zapallo italiano"
$ws.Range("H525").Value = "Sin especificar"
$ws.Range("I525").Value = "Primera"
$ws.Range("J525").Value = 400
$ws.Range("K525").Value = 11500
$ws.Range("L525").Value = 12000
$ws.Range("M525").Value = 11750
$ws.Range("N525").Value = "`$/caja 50 unidades"
$ws.Range("O525").Value = "Región de Arica y Parinacota"
$ws.Range("P525").Value = 235
$ws.Range("Q525").Value = 50
$ws.Range("R525").Value = "Hortaliza"
